$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '73.510.29'
$ws.Range("E2").Value = '  +1.49%  '

# Row 3
$ws.Range("D3").Value = '3.987.56'
$ws.Range("E3").Value = '  -1.40%  '

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = "'618.97"
$ws.Range("E5").Value = '  +14.13%  '

# Row 6
$ws.Range("D6").Value = "'168.30"
$ws.Range("E6").Value = '  +10.68%  '

# Row 7
$ws.Range("E7").Value = '  -2.25%  '

# Row 9
$ws.Range("D9").Value = "'0.760"
$ws.Range("E9").Value = '  +0.84%  '

# Row 10
$ws.Range("D10").Value = "'0.188"
$ws.Range("E10").Value = '  +8.91%  '

# Row 11
$ws.Range("D11").Value = "'55.97"
$ws.Range("E11").Value = '  +4.15%  '

# Row 12
$ws.Range("E12").Value = '  +1.63%  '

# Row 13
$ws.Range("D13").Value = "'11.19"
$ws.Range("E13").Value = '  +2.16%  '

# Row 14
$ws.Range("D14").Value = '4.624.21'
$ws.Range("E14").Value = '  -1.28%  '

# Row 15
$ws.Range("D15").Value = '3.983.58'
$ws.Range("E15").Value = '  -1.38%  '

# Row 16
$ws.Range("D16").Value = "'1.25"
$ws.Range("E16").Value = '  +3.14%  '

# Row 17
$ws.Range("D17").Value = "'14.07"
$ws.Range("E17").Value = '  -2.05%  '

# Row 18
$ws.Range("D18").Value = "'20.52"
$ws.Range("E18").Value = '  -1.12%  '

# Row 19
$ws.Range("D19").Value = '73.261.08'
$ws.Range("E19").Value = '  +1.21%  '

# Row 20
$ws.Range("E20").Value = '  -0.71%  '

# Row 21
$ws.Range("D21").Value = "'441.08"
$ws.Range("E21").Value = '  -1.95%  '

# Row 22
$ws.Range("D22").Value = "'4.88"
$ws.Range("E22").Value = '  +14.13%  '

# Row 23
$ws.Range("D23").Value = "'96.17"
$ws.Range("E23").Value = '  -1.88%  '

# Row 24
$ws.Range("D24").Value = "'3.39"
$ws.Range("E24").Value = '  -4.22%  '

# Row 25
$ws.Range("D25").Value = "'14.23"
$ws.Range("E25").Value = '  -2.81%  '

# Row 26
$ws.Range("D26").Value = "'4.10"
$ws.Range("E26").Value = '  -3.85%  '

# Row 27
$ws.Range("D27").Value = "'11.09"
$ws.Range("E27").Value = '  -1.80%  '

# Row 28
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = "'10.60"
$ws.Range("E28").Value = '  -2.06%  '

# Row 29
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").Value = "'5.97"
$ws.Range("E29").Value = '  +0.14%  '

# Row 30
$ws.Range("D30").Value = "'36.24"
$ws.Range("E30").Value = '  -2.57%  '

# Row 31
$ws.Range("D31").Value = "'7.86"
$ws.Range("E31").Value = '  -1.36%  '

# Row 32
$ws.Range("D32").Value = "'13.74"
$ws.Range("E32").Value = '  +0.84%  '

# Row 33
$ws.Range("D33").Value = "'0.0000104"
$ws.Range("E33").Value = '  +16.59%  '

# Row 34
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = '  -3.16%  '

# Row 35
$ws.Range("D35").Value = "'48.33"
$ws.Range("E35").Value = '  -1.61%  '

# Row 36
$ws.Range("D36").Value = "'71.15"
$ws.Range("E36").Value = '  +6.50%  '

# Row 37
$ws.Range("D37").Value = "'648.14"
$ws.Range("E37").Value = '  -4.86%  '

# Row 38
$ws.Range("D38").Value = "'0.431"
$ws.Range("E38").Value = '  -4.28%  '

# Row 39
$ws.Range("E39").Value = '  +0.95%  '

# Row 40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("E41").Value = '  -1.77%  '

# Row 42
$ws.Range("E42").Value = '  +0.08%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = "'0.0484"
$ws.Range("E43").Value = '  -2.37%  '

# Row 44
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = "'10.73"
$ws.Range("E44").Value = '  -4.10%  '

# Row 45
$ws.Range("D45").Value = "'3.20"
$ws.Range("E45").Value = '  -6.74%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = "'3.09"
$ws.Range("E46").Value = '  +33.24%  '

# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = "'0.149"
$ws.Range("E47").Value = '  -1.70%  '

# Row 48
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = "'0.000298"
$ws.Range("E48").Value = '  +5.76%  '

# Row 49
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = "'3.41"
$ws.Range("E49").Value = '  +3.46%  '

# Row 50
$ws.Range("D50").Value = "'2.57"
$ws.Range("E50").Value = '  -5.41%  '

# Row 51
$ws.Range("D51").Value = '2.833.04'
$ws.Range("E51").Value = '  +3.33%  '
